$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($row, $col, $val) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

# Row 2
Set-TextValue 2 4 "69.640.01"
$ws.Cells.Item(2, 5).Value = "  +0.51%  "

# Row 3
Set-TextValue 3 4 "3.687.79"
$ws.Cells.Item(3, 5).Value = "  +0.41%  "

# Row 4
Set-TextValue 4 4 "0.999"
$ws.Cells.Item(4, 5).Value = "  -0.12%  "

# Row 5
Set-TextValue 5 4 "666.54"
$ws.Cells.Item(5, 5).Value = "  -1.40%  "

# Row 6
Set-TextValue 6 4 "160.10"
$ws.Cells.Item(6, 5).Value = "  +1.04%  "

# Row 7
$ws.Cells.Item(7, 5).Value = "  +0.03%  "

# Row 8
$ws.Cells.Item(8, 5).Value = "  +1.41%  "

# Row 9
$ws.Cells.Item(9, 5).Value = "  +0.03%  "

# Row 10
Set-TextValue 10 4 "7.14"
$ws.Cells.Item(10, 5).Value = "  +3.13%  "

# Row 11
Set-TextValue 11 4 "0.442"
$ws.Cells.Item(11, 5).Value = "  +1.56%  "

# Row 12
Set-TextValue 12 4 "0.0000234"
$ws.Cells.Item(12, 5).Value = "  +1.19%  "

# Row 13
Set-TextValue 13 4 "32.95"
$ws.Cells.Item(13, 5).Value = "  +1.83%  "

# Row 14
Set-TextValue 14 4 "3.667.72"
$ws.Cells.Item(14, 5).Value = "  -0.01%  "

# Row 15
Set-TextValue 15 4 "69.606.55"

# Row 16
$ws.Cells.Item(16, 5).Value = "  +2.48%  "

# Row 17
Set-TextValue 17 4 "16.19"
$ws.Cells.Item(17, 5).Value = "  +0.88%  "

# Row 18
Set-TextValue 18 4 "6.47"
$ws.Cells.Item(18, 5).Value = "  +0.55%  "

# Row 19
Set-TextValue 19 4 "470.31"
$ws.Cells.Item(19, 5).Value = "  +0.62%  "

# Row 20
Set-TextValue 20 4 "9.77"
$ws.Cells.Item(20, 5).Value = "  -2.40%  "

# Row 21
Set-TextValue 21 4 "0.647"
$ws.Cells.Item(21, 5).Value = "  -0.24%  "

# Row 22
Set-TextValue 22 4 "79.74"
$ws.Cells.Item(22, 5).Value = "  +0.03%  "

# Row 23
Set-TextValue 23 4 "3.832.52"
$ws.Cells.Item(23, 5).Value = "  +0.36%  "

# Row 24
$ws.Cells.Item(24, 2).Value = "PEPE"
$ws.Cells.Item(24, 3).Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
Set-TextValue 24 4 "0.0000127"
$ws.Cells.Item(24, 5).Value = "  +4.60%  "

# Row 25
$ws.Cells.Item(25, 2).Value = "Dai"
$ws.Cells.Item(25, 3).Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
Set-TextValue 25 4 "1.00"
$ws.Cells.Item(25, 5).Value = "  +0.01%  "

# Row 26
Set-TextValue 26 4 "10.94"
$ws.Cells.Item(26, 5).Value = "  +0.30%  "

# Row 27
Set-TextValue 27 4 "9.06"
$ws.Cells.Item(27, 5).Value = "  -0.03%  "

# Row 28
$ws.Cells.Item(28, 5).Value = "  -0.24%  "

# Row 29
Set-TextValue 29 4 "1.71"
$ws.Cells.Item(29, 5).Value = "  -2.49%  "

# Row 30
$ws.Cells.Item(30, 5).Value = "  +0.83%  "

# Row 31
$ws.Cells.Item(31, 5).Value = "  +0.15%  "

# Row 32
Set-TextValue 32 4 "0.165"
$ws.Cells.Item(32, 5).Value = "  +2.23%  "

# Row 33
Set-TextValue 33 4 "26.75"
$ws.Cells.Item(33, 5).Value = "  -0.47%  "

# Row 34
Set-TextValue 34 4 "6.49"
$ws.Cells.Item(34, 5).Value = "  -1.76%  "

# Row 35
Set-TextValue 35 4 "3.675.34"
$ws.Cells.Item(35, 5).Value = "  +0.26%  "

# Row 36
Set-TextValue 36 4 "8.46"
$ws.Cells.Item(36, 5).Value = "  +3.44%  "

# Row 37
Set-TextValue 37 4 "6.11"
$ws.Cells.Item(37, 5).Value = "  -1.61%  "

# Row 39
$ws.Cells.Item(39, 5).Value = "  +0.61%  "

# Row 40
Set-TextValue 40 4 "0.999"
$ws.Cells.Item(40, 5).Value = "  -0.08%  "

# Row 41
Set-TextValue 41 4 "176.50"
$ws.Cells.Item(41, 5).Value = "  +1.22%  "

# Row 42
Set-TextValue 42 4 "0.0908"
$ws.Cells.Item(42, 5).Value = "  +0.86%  "

# Row 43
$ws.Cells.Item(43, 5).Value = "  -0.73%  "

# Row 44
Set-TextValue 44 4 "47.02"

# Row 45
Set-TextValue 45 4 "2.75"
$ws.Cells.Item(45, 5).Value = "  +1.79%  "

# Row 46
$ws.Cells.Item(46, 2).Value = "InjectiveProtocol"
$ws.Cells.Item(46, 3).Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextValue 46 4 "27.57"
$ws.Cells.Item(46, 5).Value = "  -2.01%  "

# Row 47
$ws.Cells.Item(47, 2).Value = "ONDO"
$ws.Cells.Item(47, 3).Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
Set-TextValue 47 4 "1.28"
$ws.Cells.Item(47, 5).Value = "  -0.67%  "

# Row 48
Set-TextValue 48 4 "0.000272"
$ws.Cells.Item(48, 5).Value = "  -2.15%  "

# Row 49
Set-TextValue 49 4 "7.86"
$ws.Cells.Item(49, 5).Value = "  +1.08%  "

# Row 50
$ws.Cells.Item(50, 5).Value = "  -0.66%  "

# Row 51
Set-TextValue 51 4 "0.264"
$ws.Cells.Item(51, 5).Value = "  -0.48%  "
